$wb = $excel.ActiveWorkbook

# --- Rename the first two project-member sheets to the actual people ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "Alwin Eizema"
$ws2.Name = "Ruud Louwerse"

# Re-point the per-sheet Print_Area so it follows the renamed sheets
$ws1.PageSetup.PrintArea = '$A$1:$D$47'
$ws2.PageSetup.PrintArea = '$A$1:$D$47'

# --- Fill in Ruud Louwerse's personal info (name / class / group) ---
$ws2.Range("B2").Value = "Ruud Louwerse"
$ws2.Range("B3").Value = "ICTM1M"
$ws2.Range("B4").Value = 2

# --- Log Ruud's timesheet entries for "Kickoff KBS" les 2 ---
$ws2.Range("A11").Value = "Kickoff KBS"
$ws2.Range("C11").Value = 120

$ws2.Range("A13").Value = "Bijeenkomst 1"
$ws2.Range("C13").Value = 30

$ws2.Range("A15").Value = "KBS les 2"
$ws2.Range("C15").Value = 120

# --- Selections / active tab: Ruud's sheet becomes the active one ---
$ws1.Activate() | Out-Null
$ws1.Range("D18").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("A15").Select() | Out-Null
